# Updated per-line active power flow (MW) results for the "380 kV" case.
# Source data: Code/Results/Cases/Case_2_195/res_line/pl_mw.xlsx, columns B:O, rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    "1.599214545649147,0.2830928052514423,0,0.2280049581199997,1.89807319322253,0.002452927055119843,0,0,0.04667749335544613,0,0.5137306350064961,0,0,3.130693982595602",
    "1.464342961756472,0.2717892734438578,0,0.2286328163377558,1.904726620059989,0.002455993654034955,0,0,0.04584481530803686,0,0.498798079695689,0,0,3.16385681041784",
    "1.381616042105634,0.2648356600484476,0,0.22909742938403,1.910015492209865,0.00245797623172033,0,0,0.04534502071383884,0,0.4897718699639455,0,0,3.186700574502481",
    "1.347927424052102,0.2619989036909089,0,0.2293066908364274,1.912473239244854,0.002458809287880393,0,0,0.04514426399664373,0,0.4861297082434248,0,0,3.196632461041872",
    "1.34233491838944,0.2615276815167817,0,0.2293426431989189,1.912899610675439,0.002458949136766388,0,0,0.0451111053231763,0,0.4855271180604888,0,0,3.198319228718887",
    "1.38116160877945,0.2647974148308947,0,0.2291001708218179,1.910047413718701,0.002457987364647778,0,0,0.04534230140360407,0,0.4897226039944655,0,0,3.186831998989106",
    "1.552694351494949,0.2791982124559524,0,0.2282050466835059,1.900117369141768,0.002453963777001599,0,0,0.04638802145516152,0,0.508552485178825,0,0,3.141612757505413",
    "1.889674101243202,0.3073253834274965,0,0.2270759022504372,1.890204852547271,0.002446860916222204,0,0,0.04852859698940293,0,0.5465982234447608,0,0,3.072682222941722",
    "2.137556349619445,0.3279131025152537,0,0.2266262434611441,1.888766877371779,0.002442117564009046,0,0,0.05015481606962879,0,0.5752236209148975,0,0,3.034150908277923",
    "2.250378369387647,0.3372605684852488,0,0.2265038068017553,1.889385585564924,0.00244006181666956,0,0,0.05090599888769276,0,0.5883903910086588,0,0,3.019268309459335",
    "2.293108108847264,0.3407974462729157,0,0.2264692182174137,1.88980316389457,0.002439297951251904,0,0,0.05119206819105671,0,0.5933969234296796,0,0,3.014014390836508",
    "2.283905238141131,0.3400358445282734,0,0.2264761442858081,1.889705075181652,0.002439461814973345,0,0,0.05113038676094561,0,0.5923177669939719,0,0,3.015128913666814",
    "2.253893652173872,0.3375516069950208,0,0.2265007253651916,1.889416265198733,0.002439998680687665,0,0,0.05092950181848011,0,0.5888018709329685,0,0,3.018828406702994",
    "2.23551147488854,0.3360295681200682,0,0.2265173145421961,1.889263237596097,0.002440329426034098,0,0,0.05080666319273064,0,0.5866509526646553,0,0,3.021144215374505",
    "2.130184182528353,0.3273018429252375,0,0.2266358945690143,1.888752077961669,0.002442253959087403,0,0,0.05010595167593834,0,0.57436603405651,0,0,3.035176873180063",
    "2.065583005952533,0.3219429023720579,0,0.2267296480772707,1.88876468572235,0.002443460679545516,0,0,0.04967898975170471,0,0.5668665560729664,0,0,3.044464010554947",
    "2.028431803405056,0.318858900926017,0,0.2267913036435161,1.888891729031855,0.002444164360831474,0,0,0.04943448791206606,0,0.5625667130128136,0,0,3.050054649937209",
    "2.015854082889803,0.3178144290732519,0,0.22681350796476,1.888955310334907,0.002444404267317513,0,0,0.04935188933753665,0,0.5611132143991853,0,0,3.051990255154493",
    "2.072459337724013,0.322513546265526,0,0.2267188679605141,1.888750943780352,0.002443331228484441,0,0,0.04972432948301986,0,0.5676634767510791,0,0,3.043449606249311",
    "2.262708622289381,0.3382813656144776,0,0.2264931859810488,1.889496119328271,0.002439840594635874,0,0,0.05098846304694149,0,0.5898340189870623,0,0,3.017731403028705",
    "2.387084433547159,0.3485701005223802,0,0.2264143114892043,1.891051587775976,0.002437644340617464,0,0,0.05182403474042729,0,0.6044434753386554,0,0,3.003149123574872",
    "2.320700004925925,0.3430803849446988,0,0.2264501401244345,1.890123558279825,0.002438808761157654,0,0,0.05137722450979254,0,0.5966352701689743,0,0,3.010727821106173",
    "2.069350582634911,0.3222555677249659,0,0.2267237174917405,1.888756783374149,0.002443389722531578,0,0,0.04970382838452636,0,0.5673031524296306,0,0,3.043907435613761",
    "1.798454278962879,0.2997292094315185,0,0.2273145162331041,1.891861243200651,0.002448698640625379,0,0,0.04794002190616453,0,0.5361869226713765,0,0,3.089208205811104"
)

$numRows = $data.Count
$numCols = 14
$vals = New-Object "object[,]" $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $fields = $data[$i].Split(",")
    for ($j = 0; $j -lt $numCols; $j++) {
        $vals[$i, $j] = [double]$fields[$j]
    }
}

$ws.Range("B2:O25").Value = $vals

